$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Info": update the summary objective/time result
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 96852329411.78194
$wsInfo.Range("B2").Value = 2.174000024795532

# ---------------------------------------------------------------------
# Sheet "Activados": Proceso column becomes 1 and the time series is
# extended from 3 points (0,140,280) to 19 points stepping by 20
# (0,20,...,360)
# ---------------------------------------------------------------------
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 0; $i -lt 19; $i++) {
    $row = 2 + $i
    $wsAct.Cells.Item($row, 1).Value = 1
    $wsAct.Cells.Item($row, 2).Value = $i * 20
}

# ---------------------------------------------------------------------
# Sheet "Operando": Proceso column changes from 4 to 1 for every data
# row (rows 2-366); Tiempo column (B) is untouched
# ---------------------------------------------------------------------
$wsOp = $wb.Worksheets.Item("Operando")
for ($row = 2; $row -le 366; $row++) {
    $wsOp.Cells.Item($row, 1).Value = 1
}

# ---------------------------------------------------------------------
# Sheet "Contaminantes": update mass / concentration results per
# contaminant (column A, the contaminant id, is untouched)
# ---------------------------------------------------------------------
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Range("B2").Value = 22460412240.00002
$wsCont.Range("C2").Value = 0.8330000000000009
$wsCont.Range("B3").Value = 1348164000.000001
$wsCont.Range("C3").Value = 0.05000000000000005
$wsCont.Range("B4").Value = 68527176119.99993
$wsCont.Range("C4").Value = 2.541499999999997
$wsCont.Range("B5").Value = 227651.7819384
$wsCont.Range("C5").Value = 0.000008443029999999999
$wsCont.Range("B6").Value = 4516349400.000105
$wsCont.Range("C6").Value = 0.1675000000000039
